$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 7605.3335
$ws.Range("I43").Value = 9822.111000000001
$ws.Range("J43").Value = 5388.5557
$ws.Range("K43").Value = 9822.111000000001
$ws.Range("L43").Value = 5388.5557
$ws.Range("M43").Value = -9753.111000000001
$ws.Range("N43").Value = -5526.5557
$ws.Range("H132").Value = 5129252
$ws.Range("I132").Value = 5129252
$ws.Range("K132").Value = 15387756
$ws.Range("M132").Value = -15385226
$ws.Range("H137").Value = 6735.4116
$ws.Range("I137").Value = 9174.375
$ws.Range("K137").Value = 27523.125
$ws.Range("M137").Value = -24973.125
$ws.Range("H138").Value = 2579.5186
$ws.Range("I138").Value = 1999.1818
$ws.Range("J138").Value = 5133
$ws.Range("K138").Value = 5997.5454
$ws.Range("L138").Value = 15399
$ws.Range("M138").Value = -857.5454
$ws.Range("N138").Value = -25679
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2682.795
$ws.Range("I61").Value = 1890.3
$ws.Range("K61").Value = 1890.3
$ws.Range("M61").Value = -1678.3
$ws.Range("H136").Value = 2682.795
$ws.Range("I136").Value = 1890.3
$ws.Range("K136").Value = 5670.9
$ws.Range("M136").Value = -3120.9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1233.75
$ws.Range("I86").Value = 1251.6
$ws.Range("J86").Value = 1144.5
$ws.Range("K86").Value = 1251.6
$ws.Range("L86").Value = 1144.5
$ws.Range("M86").Value = -128.5999999999999
$ws.Range("N86").Value = -3390.5
$ws.Range("H89").Value = 1233.75
$ws.Range("I89").Value = 1251.6
$ws.Range("J89").Value = 1144.5
$ws.Range("K89").Value = 6258
$ws.Range("L89").Value = 5722.5
$ws.Range("M89").Value = -642
$ws.Range("N89").Value = -16954.5
$ws.Range("H107").Value = 20246.555
$ws.Range("I107").Value = 23563.521
$ws.Range("J107").Value = 1174
$ws.Range("K107").Value = 23563.521
$ws.Range("L107").Value = 1174
$ws.Range("M107").Value = -21643.521
$ws.Range("N107").Value = -5014
$ws.Range("H123").Value = 79950
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 79950
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 79950
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -89750
$ws.Range("H134").Value = 1718.6
$ws.Range("I134").Value = 1500
$ws.Range("J134").Value = 3549.375
$ws.Range("K134").Value = 4500
$ws.Range("L134").Value = 10648.125
$ws.Range("M134").Value = -1965
$ws.Range("N134").Value = -15718.125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3128930
$ws.Range("I31").Value = 4169594
$ws.Range("K31").Value = 4169594
$ws.Range("M31").Value = -4169299
$ws.Range("H34").Value = 3128930
$ws.Range("I34").Value = 4169594
$ws.Range("K34").Value = 4169594
$ws.Range("M34").Value = -4169392
$ws.Range("H86").Value = 66042.914
$ws.Range("I86").Value = 146623
$ws.Range("J86").Value = 25752.875
$ws.Range("K86").Value = 146623
$ws.Range("L86").Value = 25752.875
$ws.Range("M86").Value = -145500
$ws.Range("N86").Value = -27998.875
$ws.Range("H89").Value = 66042.914
$ws.Range("I89").Value = 146623
$ws.Range("J89").Value = 25752.875
$ws.Range("K89").Value = 733115
$ws.Range("L89").Value = 128764.375
$ws.Range("M89").Value = -727499
$ws.Range("N89").Value = -139996.375
$ws.Range("H94").Value = 1395.625
$ws.Range("J94").Value = 1095.4
$ws.Range("L94").Value = 1095.4
$ws.Range("N94").Value = -1997.4
$ws.Range("H122").Value = 10765.174
$ws.Range("J122").Value = 2669.2856
$ws.Range("L122").Value = 8007.8568
$ws.Range("N122").Value = -12907.8568
$ws.Range("H132").Value = 20324.451
$ws.Range("I132").Value = 20324.451
$ws.Range("K132").Value = 60973.353
$ws.Range("M132").Value = -58443.353
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 66.95652
$ws.Range("I50").Value = 25.2
$ws.Range("J50").Value = 78.55556
$ws.Range("K50").Value = 75.59999999999999
$ws.Range("L50").Value = 235.66668
$ws.Range("M50").Value = 405.4
$ws.Range("N50").Value = -1197.66668
$ws.Range("H53").Value = 66.95652
$ws.Range("I53").Value = 25.2
$ws.Range("J53").Value = 78.55556
$ws.Range("K53").Value = 75.59999999999999
$ws.Range("L53").Value = 235.66668
$ws.Range("M53").Value = 405.4
$ws.Range("N53").Value = -1197.66668
$ws.Range("H98").Value = 2085
$ws.Range("J98").Value = 2240
$ws.Range("L98").Value = 6720
$ws.Range("N98").Value = -9716
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6542.6665
$ws.Range("I80").Value = 4672
$ws.Range("J80").Value = 11085.714
$ws.Range("K80").Value = 4672
$ws.Range("L80").Value = 11085.714
$ws.Range("M80").Value = -3674
$ws.Range("N80").Value = -13081.714
$ws.Range("H83").Value = 6542.6665
$ws.Range("I83").Value = 4672
$ws.Range("J83").Value = 11085.714
$ws.Range("K83").Value = 23360
$ws.Range("L83").Value = 55428.57
$ws.Range("M83").Value = -18368
$ws.Range("N83").Value = -65412.57
$ws.Range("H132").Value = 1440
$ws.Range("I132").Value = 1166.7
$ws.Range("K132").Value = 3500.1
$ws.Range("M132").Value = -970.1000000000004
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2504.875
$ws.Range("I136").Value = 2252.1538
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 6756.4614
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -4206.4614
$ws.Range("N136").Value = -15900
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 35000
$ws.Range("J95").Value = 35000
$ws.Range("L95").Value = 35000
$ws.Range("M95").Value = -40492
$ws.Range("H96").Value = 1895.375
$ws.Range("I96").Value = 1760.5714
$ws.Range("J96").Value = 2000.2222
$ws.Range("K96").Value = 1760.5714
$ws.Range("L96").Value = 2000.2222
$ws.Range("M96").Value = -387.5714
$ws.Range("N96").Value = -4746.2222
$ws.Range("H122").Value = 51044.184
$ws.Range("I122").Value = 78655.82000000001
$ws.Range("J122").Value = 4104.4
$ws.Range("K122").Value = 235967.46
$ws.Range("L122").Value = 12313.2
$ws.Range("M122").Value = -233517.46
$ws.Range("N122").Value = -17213.2
$ws.Range("H126").Value = 314556.38
$ws.Range("I126").Value = 2198.9092
$ws.Range("K126").Value = 6596.7276
$ws.Range("M126").Value = -4126.7276
$ws.Range("H132").Value = 2867.82
$ws.Range("I132").Value = 2933.1555
$ws.Range("J132").Value = 2279.8
$ws.Range("K132").Value = 8799.466499999999
$ws.Range("L132").Value = 6839.400000000001
$ws.Range("M132").Value = -6269.466499999999
$ws.Range("N132").Value = -11899.4
